# fix the merge error
# Player.xlsx - sheet "Property": restore the "View" (column F) flag that was
# dropped for rows 68-75 during a bad merge, and correct the boolean flags
# that ended up on the wrong rows (76-78).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Column F header is the "View" flag column (same family as Public/Private/Save).
$ws.Range("F1").Value = "View"

# Rows 68-75 lost their "View" (column F) boolean during the merge - restore TRUE.
$ws.Range("F68").Value = $true
$ws.Range("F69").Value = $true
$ws.Range("F70").Value = $true
$ws.Range("F71").Value = $true
$ws.Range("F72").Value = $true
$ws.Range("F73").Value = $true
$ws.Range("F74").Value = $true
$ws.Range("F75").Value = $true

# Rows 76-77: Private/Save flags got merged in as TRUE/TRUE/FALSE instead of
# FALSE/FALSE/TRUE - fix the Private, Save and View columns.
$ws.Range("D76").Value = $false
$ws.Range("E76").Value = $false
$ws.Range("F76").Value = $true

$ws.Range("D77").Value = $false
$ws.Range("E77").Value = $false
$ws.Range("F77").Value = $true

# Row 78: Public flag incorrectly carried over as TRUE - fix to FALSE.
$ws.Range("C78").Value = $false

# Restore the view/selection state of the sheet.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C78").Select() | Out-Null
